$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (32) had its phone number stored as text; normalize
# it to a real number (matches how the sheet re-saves existing rows).
$ws.Range("A32").Value = 71277620

# Append the new payment row (row 33).
# Phone number is numeric-looking but kept as text (same quirk as before),
# so force text storage with a leading apostrophe.
$ws.Range("A33").Value = "'71277620"
$ws.Range("B33").Value = "'"
$ws.Range("C33").Value = "Cash"
$ws.Range("D33").Value = "2025-08-18T17:10:03"
$ws.Range("E33").Value = 76
$ws.Range("F33").Value = "'"
$ws.Range("G33").Value = 76
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
